$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $orig = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $orig
}

Set-TextValue $ws.Range('D2') '26.088.84'
$ws.Range('E2').Value = '  -0.15%  '
Set-TextValue $ws.Range('D3') '1.655.35'
$ws.Range('E3').Value = '  -0.24%  '
Set-TextValue $ws.Range('D4') '1.001'
$ws.Range('E4').Value = '  -0.15%  '
Set-TextValue $ws.Range('D5') '217.63'
$ws.Range('E5').Value = '  +0.61%  '
Set-TextValue $ws.Range('D6') '0.5259'
$ws.Range('E6').Value = '  +2.22%  '
$ws.Range('E7').Value = '  -0.12%  '
Set-TextValue $ws.Range('D8') '0.2602'
$ws.Range('E8').Value = '  -1.42%  '
Set-TextValue $ws.Range('D9') '0.06342'
$ws.Range('E9').Value = '  +1.24%  '
Set-TextValue $ws.Range('D10') '20.39'
$ws.Range('E10').Value = '  -1.61%  '
Set-TextValue $ws.Range('D11') '0.07803'
$ws.Range('E11').Value = '  +0.82%  '
Set-TextValue $ws.Range('D12') '4.505'
$ws.Range('E12').Value = '  +1.62%  '
Set-TextValue $ws.Range('D13') '1.655.42'
$ws.Range('E13').Value = '  +0.13%  '
Set-TextValue $ws.Range('D14') '0.5482'
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('E15').Value = '  +1.84%  '
Set-TextValue $ws.Range('D16') '65.48'
$ws.Range('E16').Value = '  +1.13%  '
Set-TextValue $ws.Range('D17') '26.118.15'
$ws.Range('E18').Value = '  -0.16%  '
Set-TextValue $ws.Range('D19') '4.583'
$ws.Range('E19').Value = '  -0.58%  '
Set-TextValue $ws.Range('D20') '190.87'
$ws.Range('E20').Value = '  -0.42%  '
Set-TextValue $ws.Range('D21') '10.06'
$ws.Range('E21').Value = '  -0.12%  '
Set-TextValue $ws.Range('D22') '6.031'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('E23').Value = '  -0.14%  '
Set-TextValue $ws.Range('D24') '141.94'
$ws.Range('E24').Value = '  +1.63%  '
$ws.Range('E25').Value = '  +1.31%  '
Set-TextValue $ws.Range('D26') '7.245'
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  +0.15%  '
Set-TextValue $ws.Range('D29') '0.05869'
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('E30').Value = '  +0.42%  '
Set-TextValue $ws.Range('D31') '3.538'
$ws.Range('E31').Value = '  -0.77%  '
Set-TextValue $ws.Range('D32') '3.261'
$ws.Range('E32').Value = '  +0.19%  '
Set-TextValue $ws.Range('D33') '1.581'
$ws.Range('E33').Value = '  -0.85%  '
Set-TextValue $ws.Range('D34') '0.9505'
$ws.Range('E34').Value = '  -1.23%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D35') '2.411'
$ws.Range('E35').Value = '  -0.42%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D36') '2.779'
$ws.Range('E36').Value = '  +0.44%  '
Set-TextValue $ws.Range('D37') '0.5715'
$ws.Range('E37').Value = '  +1.31%  '
Set-TextValue $ws.Range('D38') '0.01618'
$ws.Range('E38').Value = '  +1.80%  '
Set-TextValue $ws.Range('D39') '5.785'
$ws.Range('E39').Value = '  -2.75%  '
Set-TextValue $ws.Range('D40') '0.8446'
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('E41').Value = '  -0.02%  '
Set-TextValue $ws.Range('D42') '103.10'
$ws.Range('E42').Value = '  +2.47%  '
Set-TextValue $ws.Range('D43') '1.026.13'
$ws.Range('E43').Value = '  +1.52%  '
Set-TextValue $ws.Range('D44') '1.799.51'
$ws.Range('E44').Value = '  -0.03%  '
Set-TextValue $ws.Range('D45') '57.11'
$ws.Range('E45').Value = '  +0.96%  '
Set-TextValue $ws.Range('D46') '1.002'
$ws.Range('E46').Value = '  -0.31%  '
Set-TextValue $ws.Range('D47') '0.4302'
$ws.Range('E47').Value = '  +2.98%  '
Set-TextValue $ws.Range('D48') '0.05150'
$ws.Range('E48').Value = '  -0.31%  '
Set-TextValue $ws.Range('D49') '1.467'
$ws.Range('E49').Value = '  +1.23%  '
Set-TextValue $ws.Range('D50') '7.806'
$ws.Range('E50').Value = '  -2.65%  '
Set-TextValue $ws.Range('D51') '0.09670'
$ws.Range('E51').Value = '  -0.19%  '

Write-Host "Applied all changes"
